$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cells with refreshed values ---
$ws.Range("G2").Value = 5047000
$ws.Range("G3").Value = 7760000
$ws.Range("G4").Value = 9320000
$ws.Range("H5").Value = 9345000
$ws.Range("H7").Value = 12112000
$ws.Range("H8").Value = 13535000
$ws.Range("G9").Value = 24802000
$ws.Range("H9").Value = 17954000
$ws.Range("G10").Value = 33913000
$ws.Range("H10").Value = 25230000
$ws.Range("G13").Value = 48777000
$ws.Range("H13").Value = 45517000
$ws.Range("G14").Value = 54410000
$ws.Range("H14").Value = 49070000
$ws.Range("G15").Value = 58807000
$ws.Range("H15").Value = 55580000
$ws.Range("G16").Value = 74830000
$ws.Range("H16").Value = 66208000
$ws.Range("C2193").Value = -0.0012101086
$ws.Range("D2193").Value = 35789800
$ws.Range("E2193").Value = 35964600
$ws.Range("F2193").Value = 35638950

# --- Append new rows 2194:2230 (data refresh through 2023-11-01) ---
$ws.Range("A2194:A2230").NumberFormat = "@"

$ws.Cells.Item(2194, 1).Value = "2023-09-26"
$ws.Cells.Item(2194, 2).Value = 35566000
$ws.Cells.Item(2194, 3).Value = 0.0021132118
$ws.Cells.Item(2194, 4).Value = 35700600
$ws.Cells.Item(2194, 5).Value = 35927400
$ws.Cells.Item(2194, 6).Value = 35665800

$ws.Cells.Item(2195, 1).Value = "2023-09-27"
$ws.Cells.Item(2195, 2).Value = 35982000
$ws.Cells.Item(2195, 3).Value = 0.0116965641
$ws.Cells.Item(2195, 4).Value = 35706600
$ws.Cells.Item(2195, 5).Value = 35949800
$ws.Cells.Item(2195, 6).Value = 35684150

$ws.Cells.Item(2196, 1).Value = "2023-09-28"
$ws.Cells.Item(2196, 2).Value = 36596000
$ws.Cells.Item(2196, 3).Value = 0.0170640876
$ws.Cells.Item(2196, 4).Value = 35833800
$ws.Cells.Item(2196, 5).Value = 36009400
$ws.Cells.Item(2196, 6).Value = 35749200

$ws.Cells.Item(2197, 1).Value = "2023-09-29"
$ws.Cells.Item(2197, 2).Value = 36444000
$ws.Cells.Item(2197, 3).Value = -0.0041534594
$ws.Cells.Item(2197, 4).Value = 36015800
$ws.Cells.Item(2197, 5).Value = 36003800
$ws.Cells.Item(2197, 6).Value = 35807400

$ws.Cells.Item(2198, 1).Value = "2023-09-30"
$ws.Cells.Item(2198, 2).Value = 36600000
$ws.Cells.Item(2198, 3).Value = 0.00428054
$ws.Cells.Item(2198, 4).Value = 36237600
$ws.Cells.Item(2198, 5).Value = 36013700
$ws.Cells.Item(2198, 6).Value = 35877600

$ws.Cells.Item(2199, 1).Value = "2023-10-01"
$ws.Cells.Item(2199, 2).Value = 37789000
$ws.Cells.Item(2199, 3).Value = 0.0324863388
$ws.Cells.Item(2199, 4).Value = 36682200
$ws.Cells.Item(2199, 5).Value = 36191400
$ws.Cells.Item(2199, 6).Value = 36050700

$ws.Cells.Item(2200, 1).Value = "2023-10-02"
$ws.Cells.Item(2200, 2).Value = 37424000
$ws.Cells.Item(2200, 3).Value = -0.009658895400000001
$ws.Cells.Item(2200, 4).Value = 36970600
$ws.Cells.Item(2200, 5).Value = 36338600
$ws.Cells.Item(2200, 6).Value = 36168050

$ws.Cells.Item(2201, 1).Value = "2023-10-03"
$ws.Cells.Item(2201, 2).Value = 37197000
$ws.Cells.Item(2201, 3).Value = -0.0060656263
$ws.Cells.Item(2201, 4).Value = 37090800
$ws.Cells.Item(2201, 5).Value = 36462300
$ws.Cells.Item(2201, 6).Value = 36253050

$ws.Cells.Item(2202, 1).Value = "2023-10-04"
$ws.Cells.Item(2202, 2).Value = 37595000
$ws.Cells.Item(2202, 3).Value = 0.0106997876
$ws.Cells.Item(2202, 4).Value = 37321000
$ws.Cells.Item(2202, 5).Value = 36668400
$ws.Cells.Item(2202, 6).Value = 36338350

$ws.Cells.Item(2203, 1).Value = "2023-10-05"
$ws.Cells.Item(2203, 2).Value = 37258000
$ws.Cells.Item(2203, 3).Value = -0.008963957999999999
$ws.Cells.Item(2203, 4).Value = 37452600
$ws.Cells.Item(2203, 5).Value = 36845100
$ws.Cells.Item(2203, 6).Value = 36404850

$ws.Cells.Item(2204, 1).Value = "2023-10-06"
$ws.Cells.Item(2204, 2).Value = 37798000
$ws.Cells.Item(2204, 3).Value = 0.0144935316
$ws.Cells.Item(2204, 4).Value = 37454400
$ws.Cells.Item(2204, 5).Value = 37068300
$ws.Cells.Item(2204, 6).Value = 36497850

$ws.Cells.Item(2205, 1).Value = "2023-10-07"
$ws.Cells.Item(2205, 2).Value = 37902000
$ws.Cells.Item(2205, 3).Value = 0.0027514683
$ws.Cells.Item(2205, 4).Value = 37550000
$ws.Cells.Item(2205, 5).Value = 37260300
$ws.Cells.Item(2205, 6).Value = 36605050

$ws.Cells.Item(2206, 1).Value = "2023-10-08"
$ws.Cells.Item(2206, 2).Value = 37929000
$ws.Cells.Item(2206, 3).Value = 0.0007123635
$ws.Cells.Item(2206, 4).Value = 37696400
$ws.Cells.Item(2206, 5).Value = 37393600
$ws.Cells.Item(2206, 6).Value = 36701500

$ws.Cells.Item(2207, 1).Value = "2023-10-09"
$ws.Cells.Item(2207, 2).Value = 37595000
$ws.Cells.Item(2207, 3).Value = -0.0088059269
$ws.Cells.Item(2207, 4).Value = 37696400
$ws.Cells.Item(2207, 5).Value = 37508700
$ws.Cells.Item(2207, 6).Value = 36756250

$ws.Cells.Item(2208, 1).Value = "2023-10-10"
$ws.Cells.Item(2208, 2).Value = 37336000
$ws.Cells.Item(2208, 3).Value = -0.006889214
$ws.Cells.Item(2208, 4).Value = 37712000
$ws.Cells.Item(2208, 5).Value = 37582300
$ws.Cells.Item(2208, 6).Value = 36798000

$ws.Cells.Item(2209, 1).Value = "2023-10-11"
$ws.Cells.Item(2209, 2).Value = 36671000
$ws.Cells.Item(2209, 3).Value = -0.0178112278
$ws.Cells.Item(2209, 4).Value = 37486600
$ws.Cells.Item(2209, 5).Value = 37470500
$ws.Cells.Item(2209, 6).Value = 36830950

$ws.Cells.Item(2210, 1).Value = "2023-10-12"
$ws.Cells.Item(2210, 2).Value = 36750000
$ws.Cells.Item(2210, 3).Value = 0.0021542909
$ws.Cells.Item(2210, 4).Value = 37256200
$ws.Cells.Item(2210, 5).Value = 37403100
$ws.Cells.Item(2210, 6).Value = 36870850

$ws.Cells.Item(2211, 1).Value = "2023-10-13"
$ws.Cells.Item(2211, 2).Value = 36968000
$ws.Cells.Item(2211, 3).Value = 0.0059319728
$ws.Cells.Item(2211, 4).Value = 37064000
$ws.Cells.Item(2211, 5).Value = 37380200
$ws.Cells.Item(2211, 6).Value = 36921250

$ws.Cells.Item(2212, 1).Value = "2023-10-14"
$ws.Cells.Item(2212, 2).Value = 36790000
$ws.Cells.Item(2212, 3).Value = -0.0048149751
$ws.Cells.Item(2212, 4).Value = 36903000
$ws.Cells.Item(2212, 5).Value = 37299700
$ws.Cells.Item(2212, 6).Value = 36984050

$ws.Cells.Item(2213, 1).Value = "2023-10-15"
$ws.Cells.Item(2213, 2).Value = 37086000
$ws.Cells.Item(2213, 3).Value = 0.0080456646
$ws.Cells.Item(2213, 4).Value = 36853000
$ws.Cells.Item(2213, 5).Value = 37282500
$ws.Cells.Item(2213, 6).Value = 37063800

$ws.Cells.Item(2214, 1).Value = "2023-10-16"
$ws.Cells.Item(2214, 2).Value = 38697000
$ws.Cells.Item(2214, 3).Value = 0.0434395729
$ws.Cells.Item(2214, 4).Value = 37258200
$ws.Cells.Item(2214, 5).Value = 37372400
$ws.Cells.Item(2214, 6).Value = 37220350

$ws.Cells.Item(2215, 1).Value = "2023-10-17"
$ws.Cells.Item(2215, 2).Value = 38652000
$ws.Cells.Item(2215, 3).Value = -0.0011628808
$ws.Cells.Item(2215, 4).Value = 37638600
$ws.Cells.Item(2215, 5).Value = 37447400
$ws.Cells.Item(2215, 6).Value = 37353850

$ws.Cells.Item(2216, 1).Value = "2023-10-18"
$ws.Cells.Item(2216, 2).Value = 38710000
$ws.Cells.Item(2216, 3).Value = 0.0015005692
$ws.Cells.Item(2216, 4).Value = 37987000
$ws.Cells.Item(2216, 5).Value = 37525500
$ws.Cells.Item(2216, 6).Value = 37459550

$ws.Cells.Item(2217, 1).Value = "2023-10-19"
$ws.Cells.Item(2217, 2).Value = 39000000
$ws.Cells.Item(2217, 3).Value = 0.0074916042
$ws.Cells.Item(2217, 4).Value = 38429000
$ws.Cells.Item(2217, 5).Value = 37666000
$ws.Cells.Item(2217, 6).Value = 37587350

$ws.Cells.Item(2218, 1).Value = "2023-10-20"
$ws.Cells.Item(2218, 2).Value = 40250000
$ws.Cells.Item(2218, 3).Value = 0.0320512821
$ws.Cells.Item(2218, 4).Value = 39061800
$ws.Cells.Item(2218, 5).Value = 37957400
$ws.Cells.Item(2218, 6).Value = 37769850

$ws.Cells.Item(2219, 1).Value = "2023-10-21"
$ws.Cells.Item(2219, 2).Value = 40400000
$ws.Cells.Item(2219, 3).Value = 0.0037267081
$ws.Cells.Item(2219, 4).Value = 39402400
$ws.Cells.Item(2219, 5).Value = 38330300
$ws.Cells.Item(2219, 6).Value = 37900400

$ws.Cells.Item(2220, 1).Value = "2023-10-22"
$ws.Cells.Item(2220, 2).Value = 40354000
$ws.Cells.Item(2220, 3).Value = -0.0011386139
$ws.Cells.Item(2220, 4).Value = 39742800
$ws.Cells.Item(2220, 5).Value = 38690700
$ws.Cells.Item(2220, 6).Value = 38046900

$ws.Cells.Item(2221, 1).Value = "2023-10-23"
$ws.Cells.Item(2221, 2).Value = 44179000
$ws.Cells.Item(2221, 3).Value = 0.0947861426
$ws.Cells.Item(2221, 4).Value = 40836600
$ws.Cells.Item(2221, 5).Value = 39411800
$ws.Cells.Item(2221, 6).Value = 38396000

$ws.Cells.Item(2222, 1).Value = "2023-10-24"
$ws.Cells.Item(2222, 2).Value = 45599000
$ws.Cells.Item(2222, 3).Value = 0.0321419679
$ws.Cells.Item(2222, 4).Value = 42156400
$ws.Cells.Item(2222, 5).Value = 40292700
$ws.Cells.Item(2222, 6).Value = 38796200

$ws.Cells.Item(2223, 1).Value = "2023-10-25"
$ws.Cells.Item(2223, 2).Value = 46484000
$ws.Cells.Item(2223, 3).Value = 0.0194083204
$ws.Cells.Item(2223, 4).Value = 43403200
$ws.Cells.Item(2223, 5).Value = 41232500
$ws.Cells.Item(2223, 6).Value = 39257500

$ws.Cells.Item(2224, 1).Value = "2023-10-26"
$ws.Cells.Item(2224, 2).Value = 46141000
$ws.Cells.Item(2224, 3).Value = -0.0073788831
$ws.Cells.Item(2224, 4).Value = 44551400
$ws.Cells.Item(2224, 5).Value = 41976900
$ws.Cells.Item(2224, 6).Value = 39674650

$ws.Cells.Item(2225, 1).Value = "2023-10-27"
$ws.Cells.Item(2225, 2).Value = 45980000
$ws.Cells.Item(2225, 3).Value = -0.0034893045
$ws.Cells.Item(2225, 4).Value = 45676600
$ws.Cells.Item(2225, 5).Value = 42709700
$ws.Cells.Item(2225, 6).Value = 40078550

$ws.Cells.Item(2226, 1).Value = "2023-10-28"
$ws.Cells.Item(2226, 2).Value = 46314000
$ws.Cells.Item(2226, 3).Value = 0.0072640278
$ws.Cells.Item(2226, 4).Value = 46103600
$ws.Cells.Item(2226, 5).Value = 43470100
$ws.Cells.Item(2226, 6).Value = 40497800

$ws.Cells.Item(2227, 1).Value = "2023-10-29"
$ws.Cells.Item(2227, 2).Value = 46825000
$ws.Cells.Item(2227, 3).Value = 0.0110333808
$ws.Cells.Item(2227, 4).Value = 46348800
$ws.Cells.Item(2227, 5).Value = 44252600
$ws.Cells.Item(2227, 6).Value = 40959300

$ws.Cells.Item(2228, 1).Value = "2023-10-30"
$ws.Cells.Item(2228, 2).Value = 46732000
$ws.Cells.Item(2228, 3).Value = -0.0019861185
$ws.Cells.Item(2228, 4).Value = 46398400
$ws.Cells.Item(2228, 5).Value = 44900800
$ws.Cells.Item(2228, 6).Value = 41429100

$ws.Cells.Item(2229, 1).Value = "2023-10-31"
$ws.Cells.Item(2229, 2).Value = 47030000
$ws.Cells.Item(2229, 3).Value = 0.0063767868
$ws.Cells.Item(2229, 4).Value = 46576200
$ws.Cells.Item(2229, 5).Value = 45563800
$ws.Cells.Item(2229, 6).Value = 41947050

$ws.Cells.Item(2230, 1).Value = "2023-11-01"
$ws.Cells.Item(2230, 2).Value = 46972000
$ws.Cells.Item(2230, 3).Value = -0.0012332554
$ws.Cells.Item(2230, 4).Value = 46777200
$ws.Cells.Item(2230, 5).Value = 46225100
$ws.Cells.Item(2230, 6).Value = 42458150

Write-Output "Update complete"